$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 51
$ws.Range("H51").Value = 1950
$ws.Range("J51").Value = 1950
$ws.Range("L51").Value = 1950
$ws.Range("N51").Value = -2918

# Row 64
$ws.Range("H64").Value = 4500
$ws.Range("I64").Value = 4500
$ws.Range("K64").Value = 4500
$ws.Range("M64").Value = -4252

# Row 67
$ws.Range("H67").Value = 4500
$ws.Range("I67").Value = 4500
$ws.Range("K67").Value = 4500
$ws.Range("M67").Value = -3642

# Row 74
$ws.Range("H74").Value = 4482
$ws.Range("I74").Value = 3866.6667
$ws.Range("J74").Value = 4943.5
$ws.Range("K74").Value = 3866.6667
$ws.Range("L74").Value = 4943.5
$ws.Range("M74").Value = -2930.6667
$ws.Range("N74").Value = -6815.5

# Row 76
$ws.Range("H76").Value = 3193.9092
$ws.Range("I76").Value = 3230.7307
$ws.Range("J76").Value = 3057.1428
$ws.Range("K76").Value = 3230.7307
$ws.Range("L76").Value = 3057.1428
$ws.Range("M76").Value = -2915.7307
$ws.Range("N76").Value = -3687.1428

# Row 77
$ws.Range("H77").Value = 4482
$ws.Range("I77").Value = 3866.6667
$ws.Range("J77").Value = 4943.5
$ws.Range("K77").Value = 19333.3335
$ws.Range("L77").Value = 24717.5
$ws.Range("M77").Value = -14653.3335
$ws.Range("N77").Value = -34077.5

# Row 79
$ws.Range("H79").Value = 3193.9092
$ws.Range("I79").Value = 3230.7307
$ws.Range("J79").Value = 3057.1428
$ws.Range("K79").Value = 3230.7307
$ws.Range("L79").Value = 3057.1428
$ws.Range("M79").Value = -2138.7307
$ws.Range("N79").Value = -5241.1428

# Row 129
$ws.Range("H129").Value = 1212.5238
$ws.Range("I129").Value = 547
$ws.Range("J129").Value = 2099.889
$ws.Range("K129").Value = 1641
$ws.Range("L129").Value = 6299.667
$ws.Range("M129").Value = 3359
$ws.Range("N129").Value = -16299.667

# Row 137
$ws.Range("H137").Value = 3558.3928
$ws.Range("I137").Value = 2375.8696
$ws.Range("K137").Value = 7127.6088
$ws.Range("M137").Value = -4577.6088

# Row 138
$ws.Range("H138").Value = 1913.3405
$ws.Range("I138").Value = 1510.4375
$ws.Range("J138").Value = 2772.8667
$ws.Range("K138").Value = 4531.3125
$ws.Range("L138").Value = 8318.6001
$ws.Range("M138").Value = 608.6875
$ws.Range("N138").Value = -18598.6001

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 583769
$ws.Range("I32").Value = 709920
$ws.Range("J32").Value = 68652.336
$ws.Range("K32").Value = 709920
$ws.Range("L32").Value = 68652.336
$ws.Range("M32").Value = -709633
$ws.Range("N32").Value = -69226.336

# Row 61
$ws.Range("H61").Value = 2443.2173
$ws.Range("I61").Value = 1524.6875
$ws.Range("J61").Value = 4542.7144
$ws.Range("K61").Value = 1524.6875
$ws.Range("L61").Value = 4542.7144
$ws.Range("M61").Value = -1312.6875
$ws.Range("N61").Value = -4966.7144

# Row 97
$ws.Range("H97").Value = 1020
$ws.Range("I97").Value = 1020
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1020
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -524
$ws.Range("N97").Value = ""

# Row 136
$ws.Range("H136").Value = 2443.2173
$ws.Range("I136").Value = 1524.6875
$ws.Range("J136").Value = 4542.7144
$ws.Range("K136").Value = 4574.0625
$ws.Range("L136").Value = 13628.1432
$ws.Range("M136").Value = -2024.0625
$ws.Range("N136").Value = -18728.1432

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1560.4
$ws.Range("I20").Value = 1504.3103
$ws.Range("K20").Value = 1504.3103
$ws.Range("M20").Value = -1257.3103

# Row 86
$ws.Range("H86").Value = 111113580
$ws.Range("I86").Value = 125002410
$ws.Range("K86").Value = 125002410
$ws.Range("M86").Value = -125001287

# Row 89
$ws.Range("H89").Value = 111113580
$ws.Range("I89").Value = 125002410
$ws.Range("K89").Value = 625012050
$ws.Range("M89").Value = -625006434

# Row 94
$ws.Range("H94").Value = 1543.6
$ws.Range("I94").Value = 1396.4615
$ws.Range("J94").Value = 2500
$ws.Range("K94").Value = 1396.4615
$ws.Range("L94").Value = 2500
$ws.Range("M94").Value = -945.4614999999999
$ws.Range("N94").Value = -3402

# Row 134
$ws.Range("H134").Value = 2702.6775
$ws.Range("I134").Value = 2391.32
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 7173.960000000001
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -4638.960000000001
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 8196.23
$ws.Range("I31").Value = 1512.1875
$ws.Range("J31").Value = 12846
$ws.Range("K31").Value = 1512.1875
$ws.Range("L31").Value = 12846
$ws.Range("M31").Value = -1217.1875
$ws.Range("N31").Value = -13436

# Row 34
$ws.Range("H34").Value = 8196.23
$ws.Range("I34").Value = 1512.1875
$ws.Range("J34").Value = 12846
$ws.Range("K34").Value = 1512.1875
$ws.Range("L34").Value = 12846
$ws.Range("M34").Value = -1310.1875
$ws.Range("N34").Value = -13250

# Row 87
$ws.Range("H87").Value = 63500.25
$ws.Range("J87").Value = 63500.25
$ws.Range("L87").Value = 63500.25
$ws.Range("N87").Value = -65872.25

# Row 90
$ws.Range("H90").Value = 63500.25
$ws.Range("J90").Value = 63500.25
$ws.Range("L90").Value = 190500.75
$ws.Range("N90").Value = -202356.75

# Row 141
$ws.Range("H141").Value = 189285.72
$ws.Range("J141").Value = 187500
$ws.Range("L141").Value = 187500
$ws.Range("N141").Value = -197860

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 707.13336
$ws.Range("I5").Value = 707.13336
$ws.Range("K5").Value = 2121.40008
$ws.Range("M5").Value = -2009.40008

# Row 34
$ws.Range("H34").Value = 11111565
$ws.Range("I34").Value = 148.57143
$ws.Range("J34").Value = 13158405
$ws.Range("K34").Value = 445.71429
$ws.Range("L34").Value = 39475215
$ws.Range("M34").Value = -361.71429
$ws.Range("N34").Value = -39475383

# Row 35
$ws.Range("H35").Value = 2925
$ws.Range("J35").Value = 5550
$ws.Range("L35").Value = 16650
$ws.Range("N35").Value = -17226

# Row 36
$ws.Range("H36").Value = 600
$ws.Range("I36").Value = 600
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 1800
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -1631
$ws.Range("N36").Value = ""

# Row 82
$ws.Range("H82").Value = 2900
$ws.Range("I82").Value = 1000
$ws.Range("J82").Value = 3111.111
$ws.Range("K82").Value = 3000
$ws.Range("L82").Value = 9333.332999999999
$ws.Range("M82").Value = -2594
$ws.Range("N82").Value = -10145.333

# Row 85
$ws.Range("H85").Value = 2900
$ws.Range("I85").Value = 1000
$ws.Range("J85").Value = 3111.111
$ws.Range("K85").Value = 3000
$ws.Range("L85").Value = 9333.332999999999
$ws.Range("M85").Value = -1596
$ws.Range("N85").Value = -12141.333

# Row 113
$ws.Range("H113").Value = 1117.381
$ws.Range("I113").Value = 619.4666999999999
$ws.Range("J113").Value = 2362.1667
$ws.Range("K113").Value = 1858.4001
$ws.Range("L113").Value = 7086.500100000001
$ws.Range("M113").Value = 311.5999000000002
$ws.Range("N113").Value = -11426.5001

# Row 118
$ws.Range("H118").Value = 2860.4546
$ws.Range("I118").Value = 963.3333
$ws.Range("K118").Value = 2889.9999
$ws.Range("M118").Value = -1646.9999

# Row 129
$ws.Range("H129").Value = 2026.2941
$ws.Range("I129").Value = 757.5
$ws.Range("J129").Value = 2416.6924
$ws.Range("K129").Value = 2272.5
$ws.Range("L129").Value = 7250.0772
$ws.Range("M129").Value = 2727.5
$ws.Range("N129").Value = -17250.0772

# Row 135
$ws.Range("H135").Value = 707.13336
$ws.Range("I135").Value = 707.13336
$ws.Range("K135").Value = 6364.20024
$ws.Range("M135").Value = -3829.20024

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5137.886
$ws.Range("I70").Value = 4914.8066
$ws.Range("K70").Value = 4914.8066
$ws.Range("M70").Value = -4644.8066

# Row 73
$ws.Range("H73").Value = 5137.886
$ws.Range("I73").Value = 4914.8066
$ws.Range("K73").Value = 4914.8066
$ws.Range("M73").Value = -3978.8066

# Row 80
$ws.Range("H80").Value = 35745960
$ws.Range("I80").Value = 45457956
$ws.Range("J80").Value = 135300
$ws.Range("K80").Value = 45457956
$ws.Range("L80").Value = 135300
$ws.Range("M80").Value = -45456958
$ws.Range("N80").Value = -137296

# Row 82
$ws.Range("H82").Value = 25999.572
$ws.Range("J82").Value = 25999.572
$ws.Range("L82").Value = 25999.572
$ws.Range("N82").Value = -26765.572

# Row 83
$ws.Range("H83").Value = 35745960
$ws.Range("I83").Value = 45457956
$ws.Range("J83").Value = 135300
$ws.Range("K83").Value = 227289780
$ws.Range("L83").Value = 676500
$ws.Range("M83").Value = -227284788
$ws.Range("N83").Value = -686484

# Row 85
$ws.Range("H85").Value = 25999.572
$ws.Range("J85").Value = 25999.572
$ws.Range("L85").Value = 25999.572
$ws.Range("N85").Value = -28651.572

# Row 94
$ws.Range("H94").Value = 57089.6
$ws.Range("J94").Value = 57089.6
$ws.Range("L94").Value = 57089.6
$ws.Range("N94").Value = -58441.6

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 4440.1
$ws.Range("I132").Value = 3501.4
$ws.Range("J132").Value = 5378.8
$ws.Range("K132").Value = 10504.2
$ws.Range("L132").Value = 16136.4
$ws.Range("M132").Value = -7974.200000000001
$ws.Range("N132").Value = -21196.4

# Row 136
$ws.Range("H136").Value = 9806481
$ws.Range("I136").Value = 2291.5833
$ws.Range("J136").Value = 33336536
$ws.Range("K136").Value = 6874.749899999999
$ws.Range("L136").Value = 100009608
$ws.Range("M136").Value = -4324.749899999999
$ws.Range("N136").Value = -100014708

$ws = $wb.Worksheets.Item("WVR")
# Row 41
$ws.Range("H41").Value = 29000
$ws.Range("J41").Value = 29000
$ws.Range("L41").Value = 29000
$ws.Range("N41").Value = -29780

# Row 81
$ws.Range("H81").Value = 6285.222
$ws.Range("I81").Value = 7942.5
$ws.Range("J81").Value = 4959.4
$ws.Range("K81").Value = 15885
$ws.Range("L81").Value = 9918.799999999999
$ws.Range("M81").Value = -14824
$ws.Range("N81").Value = -12040.8

# Row 82
$ws.Range("H82").Value = 35326.668
$ws.Range("J82").Value = 35326.668
$ws.Range("L82").Value = 35326.668
$ws.Range("N82").Value = -36092.668

# Row 84
$ws.Range("H84").Value = 6285.222
$ws.Range("I84").Value = 7942.5
$ws.Range("J84").Value = 4959.4
$ws.Range("K84").Value = 79425
$ws.Range("L84").Value = 49594
$ws.Range("M84").Value = -74121
$ws.Range("N84").Value = -60202

# Row 85
$ws.Range("H85").Value = 35326.668
$ws.Range("J85").Value = 35326.668
$ws.Range("L85").Value = 35326.668
$ws.Range("N85").Value = -37978.668

# Row 132
$ws.Range("H132").Value = 8775352
$ws.Range("J132").Value = 20836034
$ws.Range("L132").Value = 62508102
$ws.Range("N132").Value = -62513162
